$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column N (header "27-jun" + its values) into new column O,
# preserving number formats/styles, then relabel the new header to "28-jun".
$ws.Range("N1:N11").Copy($ws.Range("O1:O11"))
$ws.Range("O1").Value = "28-jun"

# Reflect the new active selection below the last used row, as in the source file.
$ws.Range("O12").Select()
